$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell A1: "الاسم الكامل" -> "الاسم"
$ws.Range("A1").Value = "الاسم"

# Update the selection to match the saved workbook state (A2)
$ws.Range("A2").Select()
